$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price cells whose new value parses as a plain number and must
# be forced to remain text (matching the source inlineStr cell type).
$dForceText = [ordered]@{
    "D5" = "300.17"
    "D6" = "96.88"
    "D7" = "0.508"
    "D9" = "0.496"
    "D11" = "0.0796"
    "D12" = "49.14"
    "D14" = "16.84"
    "D15" = "6.78"
    "D21" = "11.53"
    "D23" = "67.33"
    "D24" = "235.75"
    "D25" = "2.00"
    "D28" = "24.31"
    "D29" = "167.70"
    "D30" = "33.81"
    "D31" = "2.04"
    "D32" = "9.09"
    "D35" = "4.92"
    "D37" = "16.80"
    "D38" = "0.0692"
    "D40" = "2.81"
    "D41" = "1.75"
    "D45" = "0.0281"
    "D46" = "9.83"
    "D47" = "17.44"
    "D50" = "52.88"
    "D51" = "4.56"
}

foreach ($addr in $dForceText.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $dForceText[$addr]
}

# D-column price cells whose new value is not a valid number literal
# (multiple dots, subscript digits, etc.) so plain assignment keeps them
# as text, matching Excel's own behaviour.
$dPlain = [ordered]@{
    "D2" = "42.925.03"
    "D3" = "2.297.51"
    "D16" = "2.654.37"
    "D17" = "2.294.58"
    "D19" = "42.804.38"
    "D20" = "0.0₃0900"
    "D44" = "1.991.85"
    "D49" = "2.522.41"
}

foreach ($addr in $dPlain.Keys) {
    $ws.Range($addr).Value = $dPlain[$addr]
}

# E-column volume/percentage cells (always text: padded with spaces).
$eValues = [ordered]@{
    "E2" = "  -1.24%  "
    "E3" = "  -1.63%  "
    "E4" = "  +0.03%  "
    "E5" = "  -1.58%  "
    "E6" = "  -4.51%  "
    "E7" = "  -1.13%  "
    "E8" = "  -0.04%  "
    "E10" = "  -5.43%  "
    "E11" = "  -0.20%  "
    "E12" = "  -4.88%  "
    "E13" = "  +1.71%  "
    "E14" = "  +7.07%  "
    "E15" = "  -0.80%  "
    "E16" = "  -1.71%  "
    "E17" = "  -2.04%  "
    "E18" = "  -0.17%  "
    "E19" = "  -1.33%  "
    "E20" = "  -0.98%  "
    "E21" = "  -2.38%  "
    "E22" = "  -1.55%  "
    "E23" = "  -0.81%  "
    "E24" = "  -1.22%  "
    "E25" = "  +0.87%  "
    "E26" = "  -0.02%  "
    "E27" = "  -3.50%  "
    "E28" = "  -3.15%  "
    "E29" = "  +1.04%  "
    "E30" = "  -2.82%  "
    "E31" = "  -1.24%  "
    "E32" = "  -1.89%  "
    "E33" = "  +0.11%  "
    "E34" = "  +4.19%  "
    "E35" = "  -3.11%  "
    "E36" = "  -1.49%  "
    "E37" = "  -0.66%  "
    "E38" = "  -2.53%  "
    "E39" = "  -1.52%  "
    "E40" = "  -3.67%  "
    "E41" = "  -4.62%  "
    "E42" = "  -2.03%  "
    "E43" = "  -3.46%  "
    "E44" = "  +0.44%  "
    "E45" = "  -1.75%  "
    "E46" = "  -0.67%  "
    "E47" = "  -6.54%  "
    "E48" = "  -3.79%  "
    "E49" = "  -1.55%  "
    "E50" = "  -4.87%  "
    "E51" = "  -7.23%  "
}

foreach ($addr in $eValues.Keys) {
    $ws.Range($addr).Value = $eValues[$addr]
}
